$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 66

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/05"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = "日"
$ws.Cells.Item($row, 3).Value = 20
$ws.Cells.Item($row, 4).Value = 201
